$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to Text format so numeric-looking strings (e.g. "1.000", "0.06563")
# retain their exact literal representation instead of being parsed as numbers.
$ws.Range("D2:D51").NumberFormat = "@"

# --- Column D (Price) updates ---
$ws.Range("D2").Value = "30.396.72"
$ws.Range("D3").Value = "1.925.61"
$ws.Range("D4").Value = "1.000"
$ws.Range("D5").Value = "240.44"
$ws.Range("D6").Value = "1.001"
$ws.Range("D7").Value = "0.4754"
$ws.Range("D8").Value = "44.32"
$ws.Range("D9").Value = "0.2855"
$ws.Range("D10").Value = "0.06563"
$ws.Range("D11").Value = "19.01"
$ws.Range("D12").Value = "105.94"
$ws.Range("D13").Value = "1.920.19"
$ws.Range("D14").Value = "0.07594"
$ws.Range("D15").Value = "5.119"
$ws.Range("D16").Value = "0.6546"
$ws.Range("D17").Value = "301.64"
$ws.Range("D18").Value = "30.413.02"
$ws.Range("D19").Value = "1.000"
$ws.Range("D21").Value = "2.172.74"
$ws.Range("D22").Value = "0.000007471"
$ws.Range("D23").Value = "5.292"
$ws.Range("D25").Value = "6.264"
$ws.Range("D26").Value = "167.20"
$ws.Range("D27").Value = "9.195"
$ws.Range("D28").Value = "20.02"
$ws.Range("D29").Value = "2.014"
$ws.Range("D30").Value = "0.1113"
$ws.Range("D31").Value = "1.354"
$ws.Range("D32").Value = "4.073"
$ws.Range("D33").Value = "3.908"
$ws.Range("D34").Value = "0.04981"
$ws.Range("D35").Value = "0.7381"
$ws.Range("D36").Value = "1.144"
$ws.Range("D37").Value = "2.746"
$ws.Range("D38").Value = "0.01934"
$ws.Range("D39").Value = "2.701"
$ws.Range("D40").Value = "2.054"
$ws.Range("D41").Value = "0.8772"
$ws.Range("D42").Value = "106.80"
$ws.Range("D43").Value = "5.793"
$ws.Range("D44").Value = "69.85"
$ws.Range("D45").Value = "0.9999"
$ws.Range("D46").Value = "0.4128"
$ws.Range("D47").Value = "7.220"
$ws.Range("D48").Value = "9.254"
$ws.Range("D49").Value = "34.81"
$ws.Range("D50").Value = "0.1195"
$ws.Range("D51").Value = "0.05619"

# --- Row 21/22 swap: Coin name and Link columns ---
$ws.Range("B21").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C21").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("B22").Value = "ShibaInu"
$ws.Range("C22").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"

# --- Column E (Volume 1h) updates ---
$ws.Range("E2").Value = "  -0.04%  "
$ws.Range("E3").Value = "  +4.12%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("E5").Value = "  +3.18%  "
$ws.Range("E6").Value = "  -0.02%  "
$ws.Range("E7").Value = "  -0.10%  "
$ws.Range("E8").Value = "  +2.33%  "
$ws.Range("E9").Value = "  +4.16%  "
$ws.Range("E10").Value = "  +3.90%  "
$ws.Range("E11").Value = "  +8.36%  "
$ws.Range("E12").Value = "  +25.40%  "
$ws.Range("E13").Value = "  +3.73%  "
$ws.Range("E14").Value = "  +1.80%  "
$ws.Range("E15").Value = "  +3.53%  "
$ws.Range("E16").Value = "  +5.14%  "
$ws.Range("E17").Value = "  +23.46%  "
$ws.Range("E18").Value = "  +0.12%  "
$ws.Range("E19").Value = "  -0.04%  "
$ws.Range("E20").Value = "  +2.28%  "
$ws.Range("E21").Value = "  +3.85%  "
$ws.Range("E22").Value = "  +2.19%  "
$ws.Range("E23").Value = "  +7.92%  "
$ws.Range("E24").Value = "  -0.11%  "
$ws.Range("E25").Value = "  +6.35%  "
$ws.Range("E26").Value = "  +1.45%  "
$ws.Range("E27").Value = "  +1.43%  "
$ws.Range("E28").Value = "  +11.53%  "
$ws.Range("E29").Value = "  +7.91%  "
$ws.Range("E30").Value = "  +8.44%  "
$ws.Range("E31").Value = "  +0.69%  "
$ws.Range("E32").Value = "  +0.92%  "
$ws.Range("E33").Value = "  +2.64%  "
$ws.Range("E34").Value = "  +3.26%  "
$ws.Range("E35").Value = "  +6.30%  "
$ws.Range("E36").Value = "  +1.68%  "
$ws.Range("E37").Value = "  +1.69%  "
$ws.Range("E38").Value = "  +2.21%  "
$ws.Range("E39").Value = "  +0.80%  "
$ws.Range("E40").Value = "  +3.08%  "
$ws.Range("E41").Value = "  +0.46%  "
$ws.Range("E42").Value = "  +0.21%  "
$ws.Range("E43").Value = "  +5.36%  "
$ws.Range("E45").Value = "  -0.06%  "
$ws.Range("E46").Value = "  +1.96%  "
$ws.Range("E47").Value = "  +1.21%  "
$ws.Range("E48").Value = "  +8.32%  "
$ws.Range("E49").Value = "  +3.45%  "
$ws.Range("E50").Value = "  +0.11%  "
$ws.Range("E51").Value = "  +2.03%  "
